# Updates "想去人数" (F) / "最低票价" (G) figures across the four sheets
# of the 北京-漫展信息 workbook, matching the "output generated at 456a3b4"
# refresh commit.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 339
$ws1.Range("F6").Value = 2582
$ws1.Range("F9").Value = 526
$ws1.Range("F11").Value = 1524
$ws1.Range("F14").Value = 1496
$ws1.Range("F15").Value = 1309
$ws1.Range("F16").Value = 8
$ws1.Range("F17").Value = 524
$ws1.Range("F18").Value = 3710
$ws1.Range("F21").Value = 748
$ws1.Range("F22").Value = 13
$ws1.Range("F23").Value = 2172
$ws1.Range("F27").Value = 10
$ws1.Range("F28").Value = 1144
$ws1.Range("F29").Value = 757
$ws1.Range("F31").Value = 1007
$ws1.Range("F32").Value = 998

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 45
$ws2.Range("G3").Value = 280
$ws2.Range("F10").Value = 87
$ws2.Range("F20").Value = 471

# --- Sheet 3: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 518
$ws3.Range("F6").Value = 70

# --- Sheet 4: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 339
$ws4.Range("F5").Value = 45
$ws4.Range("G9").Value = 280
$ws4.Range("F11").Value = 518
$ws4.Range("F12").Value = 70
$ws4.Range("F13").Value = 2582
$ws4.Range("F14").Value = 2582
$ws4.Range("F20").Value = 526
$ws4.Range("F22").Value = 87
$ws4.Range("F23").Value = 1524
$ws4.Range("F26").Value = 1496
$ws4.Range("F28").Value = 1309
$ws4.Range("F29").Value = 8
$ws4.Range("F30").Value = 524
$ws4.Range("F32").Value = 3710
$ws4.Range("F35").Value = 2172
$ws4.Range("F39").Value = 10
$ws4.Range("F40").Value = 1144
$ws4.Range("F45").Value = 471
$ws4.Range("F46").Value = 757
$ws4.Range("F48").Value = 1007
$ws4.Range("F49").Value = 998
